# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# Column G ("K") on Sheet1 is recomputed from its previous "Strike#"-style
# counts into the new "s_vals" series (derived via a mean/std regen pass).
# The table runs from row 2 (r=0) through row 58 (r=56); row index is the
# value in column A, so sheet row = A-value + 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values, in sheet-row order starting at row 2.
$sVals = @(
    1,1,0,0,0,1,0,0,0,0,1,1,2,1,1,1,1,0,0,0,0,0,2,1,2,0,0,2,0,1,0,0,0,1,0,2,2,0,
    1,1,1,0,0,3,1,0,2,1,1,2,1,0,1,2,2,1,1
)

$firstRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
